# Updated cryptos list values per diff (row B/C/D/E cells)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.691.56'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +4.79%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.301.72'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.59%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.62'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.628'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.64%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '62.37'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.51%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.420'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0924'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.05%  '
$ws.Range('B11').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C11').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '2.645.06'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.68%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.104'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '16.00'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.67'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +8.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.77'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.819'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.304.51'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.65%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.495.08'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0937'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.99%  '
$ws.Range('E20').Value = '  +4.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.62'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '250.56'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.58'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +7.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.39'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.88'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '170.15'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.64%  '
$ws.Range('E28').Value = '  +1.38%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.66'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.61%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.49'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.15%  '
$ws.Range('E31').Value = '  +1.57%  '
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.06'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.76'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0660'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.68%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.48'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.56'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.67'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('E39').Value = '  +4.99%  '
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('E41').Value = '  +1.69%  '
$ws.Range('B42').Value = 'FTXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.71'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.70%  '
$ws.Range('B43').Value = 'TerraClassic'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.000219'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -12.18%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0975'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.21'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.36%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '98.70'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.00'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.04%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.475.96'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.50%  '
$ws.Range('E49').Value = '  +2.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.30'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +9.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.77'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.17%  '
